# Applies cryptos list price/volume refresh (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.108.06"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").Value = "1.870.41"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'307.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.5046"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("D8").Value = "'0.3753"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D9").Value = "'0.07159"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").Value = "'0.8902"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").Value = "'20.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "1.874.00"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "'0.07559"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "'5.329"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("D15").Value = "'89.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "'0.000008512"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.51%  "
$ws.Range("D18").Value = "'14.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.17%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "27.151.93"
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("D21").Value = "'5.091"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "2.112.68"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "'10.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").Value = "'6.495"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").Value = "'151.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").Value = "'1.841"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "'18.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "'2.096"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.58%  "
$ws.Range("D29").Value = "'112.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").Value = "'4.765"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").Value = "'4.695"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("D32").Value = "'0.08986"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").Value = "'0.05140"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").Value = "'3.098"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("D35").Value = "'0.7456"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("D36").Value = "'1.163"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.46%  "
$ws.Range("D37").Value = "'2.558"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").Value = "'0.02037"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").Value = "'3.042"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").Value = "'1.074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("D41").Value = "'0.5368"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").Value = "'6.624"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("D43").Value = "'114.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("D44").Value = "'8.479"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'0.1479"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").Value = "'0.4657"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").Value = "'10.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.49%  "
$ws.Range("D49").Value = "'1.575"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.67%  "
$ws.Range("D50").Value = "'64.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.08%  "
$ws.Range("D51").Value = "'36.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.08%  "
